# feat: add 2022-Q4 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q4" right after "总计", before "2022-Q2".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet    = $wb.Worksheets.Item("2022-Q2")
$newSheet   = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q4"

# Header row (same column layout / style as the other quarterly sheets).
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the bold/centered/bordered header style from the "总计" sheet's header.
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Data row (single fund entry for 2022-Q4).
$newSheet.Range("A2").Value = 0
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B2").Value = "'010322"
$newSheet.Range("C2").Value = "大摩新兴产业股票"
$newSheet.Range("D2").Value = "'2.06"
$newSheet.Range("E2").Value = "'93.45"
$newSheet.Range("F2").Value = "'3.58"
$newSheet.Range("G2").Value = "'0.0737"
$newSheet.Range("H2").Value = 10
# Drop the implicit "text" number-format Excel stamps on the numeric-looking
# strings above so the cells stay plain/unstyled, matching the other
# quarterly sheets (only A2 keeps the bordered/bold "index" style).
$newSheet.Range("B2:G2").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Prepend a 2022-Q4 summary row into the "总计" sheet, pushing the
#    existing 2022-Q2 / 2022-Q1 rows down by one.
# ---------------------------------------------------------------------------
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 1.24

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 1.16

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.07000000000000001

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3:A4").PasteSpecial(-4122)
